$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 (old Grating 3 / Grating 4 data rows) and clear column L (old Ih columns)
$ws.Range("A3:L4").Clear()
$ws.Range("L1:L2").Clear()

# Row 1 headers (new layout with Io columns instead of Ih columns)
$ws.Range("B1").Value = "I(uW)"
$ws.Range("D1").Value = "F1 I1(uW)"
$ws.Range("E1").Value = "F1 Iv(uW)"
$ws.Range("G1").Value = "F2 I1(uW)"
$ws.Range("H1").Value = "F2 Iv(uW)"
$ws.Range("J1").Value = "F3 I1(uW)"
$ws.Range("K1").Value = "F3 Iv(uW)"

# Row 2 data for Grating 5
$ws.Range("A2").Value = "Grating 5 "

$ws.Range("C1").Value = "F1 Io(uW)"
$ws.Range("F1").Value = "F2 Io (uW)"
$ws.Range("I1").Value = "F3 Io(uW)"

$ws.Range("B2").Value = 974.7
$ws.Range("C2").Value = 70.7
$ws.Range("D2").Value = 369.8
$ws.Range("E2").Value = 337.2
$ws.Range("F2").Value = 72.6
$ws.Range("G2").Value = 346
$ws.Range("H2").Value = 321
$ws.Range("I2").Value = 71.2
$ws.Range("J2").Value = 368
$ws.Range("K2").Value = 320.7

# Update selection to mirror the authored state
$ws.Range("G9").Select()
